# Trade #13 closed at 2026-02-17 08:14:03 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.55
$summary.Range("B6").Value = 13
$summary.Range("B9").Value = 30.77

# --- Strategy Status sheet --------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 13
$status.Range("G4").Value = 30.77

# --- New closed-trade row (row 14) data -------------------------------------
$tradeNum      = 13
$tradeDate     = "2026-02-17"
$tradeTime     = "08:13:57"
$strategy      = "MarketMaking"
$side          = "DOWN"
$entryPrice    = 0.35
$exitPrice     = 0.35
$tradeStatus   = "CLOSED"
$pnlPct        = 0
$pnlUsd        = 0
$capitalAfter  = 99.64
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = "early_exit"
$duration      = 0.13

function Add-TradeRow($ws) {
    $ws.Cells.Item(14, 1).Value = $tradeNum

    # Force text storage for the date/time strings so Excel doesn't
    # auto-coerce them into date/time serial numbers.
    $ws.Cells.Item(14, 2).NumberFormat = "@"
    $ws.Cells.Item(14, 2).Value = $tradeDate
    $ws.Cells.Item(14, 2).Style = "Normal"

    $ws.Cells.Item(14, 3).NumberFormat = "@"
    $ws.Cells.Item(14, 3).Value = $tradeTime
    $ws.Cells.Item(14, 3).Style = "Normal"

    $ws.Cells.Item(14, 4).Value = $strategy
    $ws.Cells.Item(14, 5).Value = $side
    $ws.Cells.Item(14, 6).Value = $entryPrice
    $ws.Cells.Item(14, 7).Value = $exitPrice
    $ws.Cells.Item(14, 8).Value = $tradeStatus
    $ws.Cells.Item(14, 9).Value = $pnlPct
    $ws.Cells.Item(14, 10).Value = $pnlUsd
    $ws.Cells.Item(14, 11).Value = $capitalAfter
    $ws.Cells.Item(14, 12).Value = $entrySlippage
    $ws.Cells.Item(14, 13).Value = $exitSlippage
    $ws.Cells.Item(14, 14).Value = $confidence
    $ws.Cells.Item(14, 15).Value = $entryReason
    $ws.Cells.Item(14, 16).Value = $exitReason
    $ws.Cells.Item(14, 17).Value = $duration
}

# --- All Trades sheet -------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# --- MarketMaking sheet ------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
